# Trade #50 (MarketMaking) closed - update All Trades, MarketMaking strategy
# sheet, Strategy Status summary row, and the top-level Summary sheet to
# reflect the closed trade's results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.73   # Current Capital
$summary.Range("B4").Value = 0.83      # Total P&L $
$summary.Range("B5").Value = 0.35      # Total P&L %
$summary.Range("B6").Value = 48        # Total Trades
$summary.Range("B7").Value = 27        # Winning Trades
$summary.Range("B9").Value = 56.25     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.73      # Capital
$status.Range("D6").Value = 19         # Trades
$status.Range("E6").Value = -0.08      # P&L $
$status.Range("F6").Value = -0.27      # P&L %
$status.Range("G6").Value = 57.89      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - Trade #50 row (row 51)
# Columns: G=Exit Price, H=Status, I=P&L %, J=P&L $, K=Capital After,
#          L=Exit Reason, M=Duration (min)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G51").Value = 0.37
$allTrades.Range("H51").Value = "CLOSED"
$allTrades.Range("I51").Value = 76.1905
$allTrades.Range("J51").Value = 0.16
$allTrades.Range("K51").Value = 99.73
$allTrades.Range("L51").Value = "early_exit"
$allTrades.Range("M51").Value = 0.17

# ---------------------------------------------------------------------
# MarketMaking sheet - Trade #50 row (row 22)
# Columns: G=Exit Price, H=Status, I=P&L %, J=P&L $, K=Capital After,
#          P=Exit Reason, Q=Duration (min)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G22").Value = 0.37
$mm.Range("H22").Value = "CLOSED"
$mm.Range("I22").Value = 76.1905
$mm.Range("J22").Value = 0.16
$mm.Range("K22").Value = 99.73
$mm.Range("P22").Value = "early_exit"
$mm.Range("Q22").Value = 0.17
